# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').Value = '''60.930.27'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  +3.09%  '
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').Value = '''2.613.08'
$ws.Range('D3').Style = $origStyle
$ws.Range('E4').Value = '  +0.04%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').Value = '''573.41'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.46%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').Value = '''143.63'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  +1.12%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').Value = '''2.639.94'
$ws.Range('D9').Style = $origStyle
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('E11').Value = '  +3.23%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').Value = '''0.154'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  -3.21%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').Value = '''0.369'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +6.70%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').Value = '''3.078.28'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +1.85%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').Value = '''60.926.35'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +3.05%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').Value = '''23.52'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  +4.76%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').Value = '''0.0000141'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +3.21%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').Value = '''2.626.67'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +2.07%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').Value = '''11.30'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +10.37%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').Value = '''4.67'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +3.13%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').Value = '''349.54'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +3.49%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').Value = '''7.15'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +14.59%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').Value = '''0.520'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +13.91%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').Value = '''64.12'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.51%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').Value = '''0.164'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +1.65%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').Value = '''0.995'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +0.62%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').Value = '''7.72'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +6.53%  '
$ws.Range('D29').Value = '0.0₃0797'
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('E30').Value = '  +7.50%  '
$ws.Range('E31').Value = '  -0.11%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').Value = '''6.32'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +4.15%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').Value = '''160.86'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +1.42%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').Value = '''19.54'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +2.77%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').Value = '''4.28'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +5.92%  '
$ws.Range('E36').Value = '  +10.63%  '
$ws.Range('E37').Value = '  +4.90%  '
$ws.Range('E38').Value = '  +6.01%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').Value = '''37.81'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +1.62%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').Value = '''0.857'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -1.62%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').Value = '''3.81'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +3.62%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').Value = '''298.32'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +1.83%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').Value = '''139.75'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +9.22%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').Value = '''0.0988'
$ws.Range('D44').Style = $origStyle
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').Value = '''0.606'
$ws.Range('D46').Style = $origStyle
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').Value = '''0.0549'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +2.40%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').Value = '''0.0241'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +3.98%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').Value = '''19.80'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +7.17%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').Value = '''10.71'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').Value = '''4.84'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +7.56%  '
